# Update zone_user_h master data: change zone_code from MOR to BLZ for the
# existing users, and add the new Belize (BLZ) zone users.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: zone_code MOR -> BLZ -------------------------
$ws.Range("B2").Value = "BLZ"
$ws.Range("B3").Value = "BLZ"

# --- New user rows for the BLZ zone --------------------------------------
$newUsers = @("keerthini", "akash", "avanish", "avanish-admintest", "ganesh", "mijan_32", "mijan_1")

$row = 4
foreach ($user in $newUsers) {
    $ws.Range("A$row").Value = "eng"
    $ws.Range("B$row").Value = "BLZ"
    $ws.Range("C$row").Value = $user
    $ws.Range("D$row").Value = $true
    $ws.Range("D$row").NumberFormat = '"TRUE";"TRUE";"FALSE"'
    $ws.Range("E$row").Value = "now()"
    $row = $row + 1
}
